$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: simple journal entry, date style only (reuse xf index 1) ---
[void]$ws.Range("A1").Copy($ws.Range("A11"))
$ws.Range("A11").Value = 43140
$ws.Range("B11").Value = "Avancer dans la partie connexion/inscription du site "
$ws.Range("C11").Value = "1 période"

# --- Row 12: entretien avec M. Egger, wrapped text, taller row ---
[void]$ws.Range("A1").Copy($ws.Range("A12"))
$ws.Range("A12").Value = 43140
$ws.Range("C12").Value = "1 période "
$ws.Range("B12").WrapText = $true
$ws.Range("B12").Value = "Entretien avec M. Egger pour faire un point sur ma documentation, je dois faire des modifications dans ma planification, améliorer mon MLD"
$ws.Rows.Item(12).RowHeight = 30

# --- Row 13: modification du MLD, wrapped text, even taller row ---
[void]$ws.Range("A1").Copy($ws.Range("A13"))
$ws.Range("A13").Value = 43140
$ws.Range("B13").WrapText = $true
$ws.Range("B13").Value = "Modification du MLD, j'ai ajouté une table pour le paiement des articles,`nj'ai aussi enlevé deux tables et ajouté deux nouveaux champs, dans la table article, qui sont `ntaille et couleur "
$ws.Range("C13").Value = "1 période"
$ws.Rows.Item(13).RowHeight = 60

# --- Update the active selection to reflect the new bottom of the log ---
[void]$ws.Range("C14").Select()

Write-Host "applied journal updates"
